$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 34 ("TV Guide") to hold the new
# "touch support" feature entry, shifting everything below it down by one.
$ws.Rows(33).Insert() | Out-Null

# New row content: Feature/bug | Status | Comments/Issues
$ws.Range("A33").Value = "Some changes to improve touch support on small screens (it's not just for mice!)"
$ws.Range("B33").Value = "Done"
$ws.Range("C33").Value = "Now you don't have to use Estouchy."

# Match the look of the other feature rows (unbolded, wrapped, top-aligned)
# and give it the same taller row height used by similar two-line rows.
$ws.Range("A33:C33").Font.Bold = $false
$ws.Rows(33).RowHeight = 28.8

# Restore the cursor position recorded in the saved workbook.
$ws.Range("C34").Select() | Out-Null
